# Update "想去人数" (interested-people count) figures on the 展览 (Exhibition)
# and 全部类型 (All types) sheets, matching the refreshed data snapshot.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 116
    $ws.Range("F5").Value = 3053
    $ws.Range("F6").Value = 309
}
